$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 230, shifting existing rows 230:264 down to 231:265
$ws.Rows.Item(230).Insert()

# Populate the newly inserted row 230 with its data
$ws.Cells.Item(230, 1).Value = 10
$ws.Cells.Item(230, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(230, 3).Value = "La Araucanía"
$ws.Cells.Item(230, 4).Value = 44522
$ws.Cells.Item(230, 5).Value = 9
$ws.Cells.Item(230, 6).Value = 100112040
$ws.Cells.Item(230, 7).Value = "Cilantro"
$ws.Cells.Item(230, 8).Value = "Sin especificar"
$ws.Cells.Item(230, 9).Value = "Primera"
$ws.Cells.Item(230, 10).Value = 40
$ws.Cells.Item(230, 11).Value = 6000
$ws.Cells.Item(230, 12).Value = 6000
$ws.Cells.Item(230, 13).Value = 6000
$ws.Cells.Item(230, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(230, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(230, 16).Value = 3000
$ws.Cells.Item(230, 17).Value = 2
$ws.Cells.Item(230, 18).Value = "Hortaliza"

# Match the date style used by column D in other rows
$ws.Cells.Item(230, 4).NumberFormat = $ws.Cells.Item(231, 4).NumberFormat
